# Generate Report for Archive
# Updates the localization-status workbook: the "Ready for handoff" status
# text becomes "In Translation" everywhere it's used, and the now-narrower
# status column is resized on each sheet that shows it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text lives in E2/F2, columns E & F hold it ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Columns.Item(5).ColumnWidth = 13.4101845877511
$ws.Columns.Item(6).ColumnWidth = 13.4101845877511

# --- zh-cn sheet: status text in C2, column C holds it ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "In Translation"
$ws.Columns.Item(3).ColumnWidth = 13.4101845877511

# --- de-de sheet: status text in C2, column C holds it ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "In Translation"
$ws.Columns.Item(3).ColumnWidth = 13.4101845877511
